$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D7 from "Open" to "Closed"
$ws.Range("D7").Value = "Closed"

# Append new rows 18 through 58
$newRows = @(
    @('Jyoti Goel', 'R90MZVZ2', 'jyoti@awlindia.com', 'Closed'),
    @('Mamta Kumari', '7XMHVP2', 'mamta.kumari@awlindia.com', 'Closed'),
    @('Sanjay Singh', '48SJ5S2', 'sanjay.singh@awlindia.com', 'Closed'),
    @('Sanjay Singh', '48SJ5S2', 'sanjay.singh@awlindia.com', 'Closed'),
    @('Sanjay Singh', '48SJ5S2', 'sanjay.singh@awlindia.com', 'Closed'),
    @('Vikas Tiwari', 'CNU406B4BL', 'vikas.tiwari@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('Lakshmi Narayan Sharma', '8ZX7HR2', 'lnsharma@awlindia.com', 'Closed'),
    @('Deepak Behera', 'DKG1T32', 'deepak@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'L9HK853', 'support@awlindia.com', 'Closed'),
    @('Arun Prakash M', '33XNGH2', 'arun.prakash@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('Devender Kumar', '5CD8135YDV', '', 'Closed'),
    @('Subham Singh', '5CD60205YD', 'subham@awlindia.com', 'Closed'),
    @('Devender Kumar', '3XNJ732', '', 'Closed'),
    @('Charu Chauhan', '5CD8135YH9', 'charu.chauhan@awlindia.com', 'Closed'),
    @('Dinesh Kumar', '9L4K5S2', 'dinesh@awlindia.com', 'Closed'),
    @('Lakshmi Narayan Sharma', '8ZX7HR2', 'lnsharma@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('Subham Singh', '5CD60205YD', 'subham@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('Rupesh Kumar', '1BN3TC2', 'Rupesh.kumar@ilogsolution.com', 'Closed'),
    @('Charu Chauhan', '5CD8135YH9', 'charu.chauhan@awlindia.com', 'Closed'),
    @('Lakshmi Narayan Sharma', '8ZX7HR2', 'lnsharma@awlindia.com', 'Closed'),
    @('Devender Kumar', '3XNJ732', '', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('Samridhi Pandey', 'FG45F12', 'samridhi@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed'),
    @('Jyoti Goel', 'R90MZVZ2', 'jyoti@awlindia.com', 'Closed'),
    @('Neeraj Kumar', 'CNU3319V79', 'neeraj@awlindia.com', 'Closed'),
    @('Akshay Kumar', 'CNU4369GOV', 'akshay.kumar@awlindia.com', 'Closed'),
    @('Vikas Tiwari', 'CNU406B4BL', 'vikas.tiwari@awlindia.com', 'Closed'),
    @('IT SUPPORT', '059881F81220034AE', 'support@awlindia.com', 'Closed'),
    @('Akshay Kumar', 'CNU4369GOV', 'akshay.kumar@awlindia.com', 'Closed'),
    @('Ankush Rastogi', 'D7CN6S2', 'ankush@awlindia.com', 'Closed'),
    @('Nadeem Khan', '4T0WKV1', 'nadeem@awlindia.com', 'Closed'),
    @('IT SUPPORT', 'Other', 'support@awlindia.com', 'Closed')
)

$startRow = 18
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne '') {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
}

